$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the sheet (ExtensionSchemes_O1234567890555 -> Extensions_O1234567890555)
$ws.Name = "Extensions_O1234567890555"

# 2. Fix the shared-string text used by F2 ("Test extensionscheme55" -> "Test extension55")
$ws.Range("F2").Value = "Test extension55"

# 3. Swap the cell formatting between C2 and F2
#    (use a scratch cell so neither source is overwritten before it's copied)
$ws.Range("C2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("Z1").Copy()
$ws.Range("F2").PasteSpecial(-4122)

$ws.Range("Z1").Clear()
